$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7,8).Value = 8500  # H7: was 7571.2856
$ws.Cells.Item(7,9).Value = 0  # I7: was 4499
$ws.Cells.Item(7,10).Value = 8500  # J7: was 8083.3335
$ws.Cells.Item(7,11).Value = 0  # K7: was 4499
$ws.Cells.Item(7,12).ClearContents()  # L7: was 8083.3335
$ws.Cells.Item(7,13).Value = 8500  # M7: was -4387
$ws.Cells.Item(7,14).Value = -8724  # N7: was -8307.333500000001

$ws.Cells.Item(8,8).Value = 520.9091  # H8: was 166.125
$ws.Cells.Item(8,9).Value = 163.33333  # I8: was 166.125
$ws.Cells.Item(8,10).Value = 950  # J8: was 0
$ws.Cells.Item(8,11).Value = 489.99999  # K8: was 498.375
$ws.Cells.Item(8,12).Value = 2850  # L8: was 0
$ws.Cells.Item(8,13).Value = -350.99999  # M8: was -359.375
$ws.Cells.Item(8,14).Value = -3128  # N8: was None

$ws.Cells.Item(14,8).Value = 8500  # H14: was 7571.2856
$ws.Cells.Item(14,9).Value = 0  # I14: was 4499
$ws.Cells.Item(14,10).Value = 8500  # J14: was 8083.3335
$ws.Cells.Item(14,11).Value = 0  # K14: was 4499
$ws.Cells.Item(14,12).ClearContents()  # L14: was 8083.3335
$ws.Cells.Item(14,13).Value = 8500  # M14: was -4308
$ws.Cells.Item(14,14).Value = -8882  # N14: was -8465.333500000001

$ws.Cells.Item(28,8).Value = 621.28125  # H28: was 638.0968
$ws.Cells.Item(28,10).Value = 1176.0769  # J28: was 1265.75
$ws.Cells.Item(28,12).Value = 1176.0769  # L28: was 1265.75
$ws.Cells.Item(28,14).Value = -2146.0769  # N28: was -2235.75

$ws.Cells.Item(34,8).Value = 3264.3  # H34: was 4269.5713
$ws.Cells.Item(34,9).Value = 1404.7778  # I34: was 1647.8334
$ws.Cells.Item(34,11).Value = 1404.7778  # K34: was 1647.8334
$ws.Cells.Item(34,13).Value = -1201.7778  # M34: was -1444.8334

$ws.Cells.Item(36,8).Value = 3264.3  # H36: was 4269.5713
$ws.Cells.Item(36,9).Value = 1404.7778  # I36: was 1647.8334
$ws.Cells.Item(36,11).Value = 1404.7778  # K36: was 1647.8334
$ws.Cells.Item(36,13).Value = -689.7778000000001  # M36: was -932.8334

$ws.Cells.Item(121,8).Value = 797.5  # H121: was 1015
$ws.Cells.Item(121,10).Value = 990  # J121: was 1296.6666
$ws.Cells.Item(121,12).Value = 2970  # L121: was 3889.9998
$ws.Cells.Item(121,14).Value = -6464  # N121: was -7383.9998

$ws.Cells.Item(135,8).Value = 494.76666  # H135: was 675
$ws.Cells.Item(135,9).Value = 494.76666  # I135: was 488.7097
$ws.Cells.Item(135,10).Value = 0  # J135: was 1500
$ws.Cells.Item(135,11).Value = 4452.89994  # K135: was 4398.3873
$ws.Cells.Item(135,12).Value = 0  # L135: was 13500
$ws.Cells.Item(135,13).ClearContents()  # M135: was -1863.3873
$ws.Cells.Item(135,14).Value = -1917.89994  # N135: was -18570

$ws.Cells.Item(138,8).Value = 2108.977  # H138: was 2123.519
$ws.Cells.Item(138,9).Value = 1004.3261  # I138: was 1022.1111
$ws.Cells.Item(138,10).Value = 3348.3416  # J138: was 3581.2646
$ws.Cells.Item(138,11).Value = 3012.9783  # K138: was 3066.3333
$ws.Cells.Item(138,12).Value = 10045.0248  # L138: was 10743.7938
$ws.Cells.Item(138,13).Value = 2127.0217  # M138: was 2073.6667
$ws.Cells.Item(138,14).Value = -20325.0248  # N138: was -21023.7938

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(46,8).Value = 2272.4443  # H46: was 6488
$ws.Cells.Item(46,9).Value = 1433.3334  # I46: was 7000
$ws.Cells.Item(46,10).Value = 2692  # J46: was 6317.3335
$ws.Cells.Item(46,11).Value = 1433.3334  # K46: was 7000
$ws.Cells.Item(46,12).Value = 2692  # L46: was 6317.3335
$ws.Cells.Item(46,13).Value = -1114.3334  # M46: was -6681
$ws.Cells.Item(46,14).Value = -3330  # N46: was -6955.3335

$ws.Cells.Item(88,8).Value = 2798  # H88: was 2908
$ws.Cells.Item(88,9).Value = 2700  # I88: was 3100
$ws.Cells.Item(88,10).Value = 2863.3333  # J88: was 2825.7144
$ws.Cells.Item(88,11).Value = 2700  # K88: was 3100
$ws.Cells.Item(88,12).Value = 2863.3333  # L88: was 2825.7144
$ws.Cells.Item(88,13).Value = -2294  # M88: was -2694
$ws.Cells.Item(88,14).Value = -3675.3333  # N88: was -3637.7144

$ws.Cells.Item(91,8).Value = 2798  # H91: was 2908
$ws.Cells.Item(91,9).Value = 2700  # I91: was 3100
$ws.Cells.Item(91,10).Value = 2863.3333  # J91: was 2825.7144
$ws.Cells.Item(91,11).Value = 2700  # K91: was 3100
$ws.Cells.Item(91,12).Value = 2863.3333  # L91: was 2825.7144
$ws.Cells.Item(91,13).Value = -1296  # M91: was -1696
$ws.Cells.Item(91,14).Value = -5671.3333  # N91: was -5633.7144

$ws.Cells.Item(98,8).Value = 19611  # H98: was 22222
$ws.Cells.Item(98,10).Value = 19611  # J98: was 22222
$ws.Cells.Item(98,12).Value = 19611  # L98: was 22222
$ws.Cells.Item(98,14).Value = -25601  # N98: was -28212

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86,8).Value = 63155.945  # H86: was 42380.832
$ws.Cells.Item(86,9).Value = 80178.57000000001  # I86: was 54248.61
$ws.Cells.Item(86,10).Value = 3576.75  # J86: was 3386.7144
$ws.Cells.Item(86,11).Value = 80178.57000000001  # K86: was 54248.61
$ws.Cells.Item(86,12).Value = 3576.75  # L86: was 3386.7144
$ws.Cells.Item(86,13).Value = -79055.57000000001  # M86: was -53125.61
$ws.Cells.Item(86,14).Value = -5822.75  # N86: was -5632.7144

$ws.Cells.Item(89,8).Value = 63155.945  # H89: was 42380.832
$ws.Cells.Item(89,9).Value = 80178.57000000001  # I89: was 54248.61
$ws.Cells.Item(89,10).Value = 3576.75  # J89: was 3386.7144
$ws.Cells.Item(89,11).Value = 400892.85  # K89: was 271243.05
$ws.Cells.Item(89,12).Value = 17883.75  # L89: was 16933.572
$ws.Cells.Item(89,13).Value = -395276.85  # M89: was -265627.05
$ws.Cells.Item(89,14).Value = -29115.75  # N89: was -28165.572

$ws.Cells.Item(134,8).Value = 2421.4062  # H134: was 2427.1875
$ws.Cells.Item(134,9).Value = 2352.3667  # I134: was 2358.5334
$ws.Cells.Item(134,11).Value = 7057.1001  # K134: was 7075.600199999999
$ws.Cells.Item(134,13).Value = -4522.1001  # M134: was -4540.600199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62,8).Value = 2581.818  # H62: was 2495
$ws.Cells.Item(62,9).Value = 2500  # I62: was 2027.5
$ws.Cells.Item(62,10).Value = 2590  # J62: was 2628.5715
$ws.Cells.Item(62,11).Value = 2500  # K62: was 2027.5
$ws.Cells.Item(62,12).Value = 2590  # L62: was 2628.5715
$ws.Cells.Item(62,13).Value = -1876  # M62: was -1403.5
$ws.Cells.Item(62,14).Value = -3838  # N62: was -3876.5715

$ws.Cells.Item(65,8).Value = 2581.818  # H65: was 2495
$ws.Cells.Item(65,9).Value = 2500  # I65: was 2027.5
$ws.Cells.Item(65,10).Value = 2590  # J65: was 2628.5715
$ws.Cells.Item(65,11).Value = 12500  # K65: was 10137.5
$ws.Cells.Item(65,12).Value = 12950  # L65: was 13142.8575
$ws.Cells.Item(65,13).Value = -9380  # M65: was -7017.5
$ws.Cells.Item(65,14).Value = -19190  # N65: was -19382.8575

$ws.Cells.Item(96,8).Value = 19965.75  # H96: was 20081.334
$ws.Cells.Item(96,10).Value = 19965.75  # J96: was 20081.334
$ws.Cells.Item(96,12).Value = 19965.75  # L96: was 20081.334
$ws.Cells.Item(96,14).Value = -25457.75  # N96: was -25573.334

$ws.Cells.Item(119,8).Value = 39989.5  # H119: was 39995.5
$ws.Cells.Item(119,10).Value = 39989.5  # J119: was 39995.5
$ws.Cells.Item(119,12).Value = 39989.5  # L119: was 39995.5
$ws.Cells.Item(119,14).Value = -49665.5  # N119: was -49671.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(105,8).Value = 133132.25  # H105: was 119145.445
$ws.Cells.Item(105,10).Value = 133132.25  # J105: was 119145.445
$ws.Cells.Item(105,12).Value = 399396.75  # L105: was 357436.335
$ws.Cells.Item(105,14).Value = -404638.75  # N105: was -362678.335

$ws.Cells.Item(131,8).Value = 823.4  # H131: was 823.87
$ws.Cells.Item(131,10).Value = 857.3913  # J131: was 857.90216
$ws.Cells.Item(131,12).Value = 2572.1739  # L131: was 2573.70648
$ws.Cells.Item(131,14).Value = -12652.1739  # N131: was -12653.70648

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22,8).Value = 1179.25  # H22: was 1104
$ws.Cells.Item(22,9).Value = 1069.3334  # I22: was 1104
$ws.Cells.Item(22,10).Value = 1509  # J22: was 0
$ws.Cells.Item(22,11).Value = 1069.3334  # K22: was 1104
$ws.Cells.Item(22,12).Value = 1509  # L22: was 0
$ws.Cells.Item(22,13).Value = -540.3334  # M22: was -575
$ws.Cells.Item(22,14).Value = -2567  # N22: was None

$ws.Cells.Item(70,8).Value = 187845.1  # H70: was 93385.87
$ws.Cells.Item(70,9).Value = 337631.16  # I70: was 139131.27
$ws.Cells.Item(70,10).Value = 8101.8  # J70: was 7613.25
$ws.Cells.Item(70,11).Value = 337631.16  # K70: was 139131.27
$ws.Cells.Item(70,12).Value = 8101.8  # L70: was 7613.25
$ws.Cells.Item(70,13).Value = -337361.16  # M70: was -138861.27
$ws.Cells.Item(70,14).Value = -8641.799999999999  # N70: was -8153.25

$ws.Cells.Item(73,8).Value = 187845.1  # H73: was 93385.87
$ws.Cells.Item(73,9).Value = 337631.16  # I73: was 139131.27
$ws.Cells.Item(73,10).Value = 8101.8  # J73: was 7613.25
$ws.Cells.Item(73,11).Value = 337631.16  # K73: was 139131.27
$ws.Cells.Item(73,12).Value = 8101.8  # L73: was 7613.25
$ws.Cells.Item(73,13).Value = -336695.16  # M73: was -138195.27
$ws.Cells.Item(73,14).Value = -9973.799999999999  # N73: was -9485.25

$ws.Cells.Item(95,8).Value = 18844  # H95: was 19800
$ws.Cells.Item(95,10).Value = 18844  # J95: was 19800
$ws.Cells.Item(95,12).Value = 18844  # L95: was 19800
$ws.Cells.Item(95,14).Value = -24336  # N95: was -25292

$ws.Cells.Item(120,8).Value = 35215.332  # H120: was 35462.6
$ws.Cells.Item(120,10).Value = 35215.332  # J120: was 35462.6
$ws.Cells.Item(120,12).Value = 35215.332  # L120: was 35462.6
$ws.Cells.Item(120,14).Value = -44891.332  # N120: was -45138.6

$ws.Cells.Item(126,8).Value = 3925714.8  # H126: was 5044.8335
$ws.Cells.Item(126,9).Value = 4652  # I126: was 5222.5713
$ws.Cells.Item(126,10).Value = 6539756.5  # J126: was 4796
$ws.Cells.Item(126,11).Value = 13956  # K126: was 15667.7139
$ws.Cells.Item(126,12).Value = 19619269.5  # L126: was 14388
$ws.Cells.Item(126,13).Value = -11486  # M126: was -13197.7139
$ws.Cells.Item(126,14).Value = -19624209.5  # N126: was -19328

$ws.Cells.Item(135,8).Value = 28615.264  # H135: was 28673.5
$ws.Cells.Item(135,10).Value = 28615.264  # J135: was 28673.5
$ws.Cells.Item(135,12).Value = 28615.264  # L135: was 28673.5
$ws.Cells.Item(135,14).Value = -38755.264  # N135: was -38813.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 2883.318  # H7: was 3194.389
$ws.Cells.Item(7,9).Value = 1794.4615  # I7: was 1919.4
$ws.Cells.Item(7,10).Value = 4456.1113  # J7: was 4788.125
$ws.Cells.Item(7,11).Value = 1794.4615  # K7: was 1919.4
$ws.Cells.Item(7,12).Value = 4456.1113  # L7: was 4788.125
$ws.Cells.Item(7,13).Value = -1682.4615  # M7: was -1807.4
$ws.Cells.Item(7,14).Value = -4680.1113  # N7: was -5012.125

$ws.Cells.Item(16,8).Value = 63406.75  # H16: was 32785.484
$ws.Cells.Item(16,9).Value = 91509.63  # I16: was 38786.46
$ws.Cells.Item(16,11).Value = 91509.63  # K16: was 38786.46
$ws.Cells.Item(16,13).Value = -91339.63  # M16: was -38616.46

$ws.Cells.Item(126,8).Value = 2883.318  # H126: was 3194.389
$ws.Cells.Item(126,9).Value = 1794.4615  # I126: was 1919.4
$ws.Cells.Item(126,10).Value = 4456.1113  # J126: was 4788.125
$ws.Cells.Item(126,11).Value = 5383.3845  # K126: was 5758.200000000001
$ws.Cells.Item(126,12).Value = 13368.3339  # L126: was 14364.375
$ws.Cells.Item(126,13).Value = -2913.3845  # M126: was -3288.200000000001
$ws.Cells.Item(126,14).Value = -18308.3339  # N126: was -19304.375

$ws.Cells.Item(132,8).Value = 2962.6667  # H132: was 3004.439
$ws.Cells.Item(132,9).Value = 3066.7812  # I132: was 3211.2334
$ws.Cells.Item(132,10).Value = 2629.5  # J132: was 2440.4546
$ws.Cells.Item(132,11).Value = 9200.3436  # K132: was 9633.700199999999
$ws.Cells.Item(132,12).Value = 7888.5  # L132: was 7321.3638
$ws.Cells.Item(132,13).Value = -6670.3436  # M132: was -7103.700199999999
$ws.Cells.Item(132,14).Value = -12948.5  # N132: was -12381.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(112,8).Value = 26387  # H112: was 29500
$ws.Cells.Item(112,10).Value = 26387  # J112: was 29500
$ws.Cells.Item(112,12).Value = 26387  # L112: was 29500
$ws.Cells.Item(112,14).Value = -29341  # N112: was -32454

$ws.Cells.Item(135,8).Value = 48249.5  # H135: was 48998
$ws.Cells.Item(135,10).Value = 48249.5  # J135: was 48998
$ws.Cells.Item(135,12).Value = 48249.5  # L135: was 48998
$ws.Cells.Item(135,14).Value = -58389.5  # N135: was -59138

$ws.Cells.Item(136,8).Value = 785.02325  # H136: was 757
$ws.Cells.Item(136,9).Value = 403  # I136: was 385.86206
$ws.Cells.Item(136,11).Value = 1209  # K136: was 1157.58618
$ws.Cells.Item(136,13).Value = 1341  # M136: was 1392.41382
